$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Swap the Approved/Rejected status (column I) and ReasonToReject (column J)
# between row 24 and row 25.
$ws.Range("I24").Value = "Rejected"
$ws.Range("J24").Value = "Nil"

$ws.Range("I25").Value = "Approved"
$ws.Range("J25").ClearContents()

# Update the selected range shown in the sheet view
$ws.Activate()
[void]$ws.Range("J22").Select()
